# "Correct the weird bug": the HR report sheet was missing a "Date of Last
# Update" column and a "TicketID" echo column before the categorical
# breakdown columns. Column I's header was mistakenly showing the category
# text while actually holding the second timestamp value. Fix the header
# text and insert the two missing columns with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I1 actually holds a timestamp (I2), not a category - relabel the header.
$ws.Range("I1").Value = "Date of Last Update"

# Insert two missing columns right after I (before the old J/K), shifting
# the existing "Human Resource / Personnel Issues" / "Location of Issue"
# columns two places to the right.
$ws.Range("J1:K1").EntireColumn.Insert()

# New column J: repeats the TicketID for this row.
$ws.Range("J1").Value = "TicketID"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "246"
$ws.Range("J2").Style = "Normal"

# New column K: the HR process/personnel category detail.
$ws.Range("K1").Value = "HR Process & Personnel Related issues"
$ws.Range("K2").Value = "Other HR-related issue"
